{"js": "// Update the Renaissance ZGC benchmark stats table:\n//   - rows 0-2 (\"99.99\", \"0.03\", \"334\") become placeholder \"0M\" values\n//   - 10 new rows (per-iteration stats) are inserted right after them\n//   - the old multi-run \"summary\" rows at the end of the table (which\n//     held the real values, tab-separated, in a single run) are\n//     collapsed down to the plain values \"99.99\", \"0.03\", \"334\" that\n//     used to live in rows 0-2\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  return \"no table found\";\n}\n\n// Step 1: the first three rows turn into \"0M\" placeholders.\nconst cellA = table.getCellOrNullObject(0, 0);\nconst cellB = table.getCellOrNullObject(1, 0);\nconst cellC = table.getCellOrNullObject(2, 0);\nawait context.sync();\n\ncellA.value = \"0M\";\ncellB.value = \"0M\";\ncellC.value = \"0M\";\nawait context.sync();\n\n// Step 2: insert the 10 new data rows right after (old) row index 2.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst anchorRow = rows.items[2];\nconst newRowValues = [\n  \"768\",\n  \"0.00003\",\n  \"0.00026\",\n  \"0.00004\",\n  \"0.00001\",\n  \"0.00004\",\n  \"0.00005\",\n  \"0.00005\",\n  \"0.03463\",\n  \"100.0\",\n];\nanchorRow.insertRows(\n  \"After\",\n  newRowValues.length,\n  newRowValues.map((v) => [v])\n);\nawait context.sync();\n\n// Step 3: the final three rows (previously multi-run \"665\\t...\\t100.0\"\n// style summaries) become plain single values.\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst lastIndex = table.rowCount - 1;\nconst tailCell0 = table.getCellOrNullObject(lastIndex - 2, 0);\nconst tailCell1 = table.getCellOrNullObject(lastIndex - 1, 0);\nconst tailCell2 = table.getCellOrNullObject(lastIndex, 0);\nawait context.sync();\n\ntailCell0.value = \"99.99\";\ntailCell1.value = \"0.03\";\ntailCell2.value = \"334\";\nawait context.sync();\n", "ps1": "# Update the Renaissance ZGC benchmark stats table:\n#   - rows 1-3 (1-indexed: \"99.99\", \"0.03\", \"334\") become placeholder \"0M\" values\n#   - 10 new rows (per-iteration stats) are inserted right after them\n#   - the old multi-run \"summary\" rows at the end of the table (which\n#     held the real values, tab-separated, in a single run) are\n#     collapsed down to the plain values \"99.99\", \"0.03\", \"334\" that\n#     used to live in rows 1-3\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Step 1: the first three rows turn into \"0M\" placeholders.\n$table.Cell(1, 1).Range.Text = \"0M\"\n$table.Cell(2, 1).Range.Text = \"0M\"\n$table.Cell(3, 1).Range.Text = \"0M\"\n\n# Step 2: insert the 10 new data rows right after (1-indexed) row 3,\n# i.e. right before the row that currently holds \"0\" (row 4).\n$newRowValues = @(\"768\", \"0.00003\", \"0.00026\", \"0.00004\", \"0.00001\", \"0.00004\", \"0.00005\", \"0.00005\", \"0.03463\", \"100.0\")\n$beforeRow = $table.Rows.Item(4)\nforeach ($v in $newRowValues) {\n    $table.Rows.Add($beforeRow) | Out-Null\n}\nfor ($i = 0; $i -lt $newRowValues.Length; $i++) {\n    $table.Cell(4 + $i, 1).Range.Text = $newRowValues[$i]\n}\n\n# Step 3: the final three rows (previously multi-run \"665`t...`t100.0\"\n# style summaries) become plain single values.\n$rowCount = $table.Rows.Count\n$table.Cell($rowCount - 2, 1).Range.Text = \"99.99\"\n$table.Cell($rowCount - 1, 1).Range.Text = \"0.03\"\n$table.Cell($rowCount, 1).Range.Text = \"334\"\n"}
